$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.757.58'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.599.19'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '''211.89'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '''0.512'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').Value = '''0.0618'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').Value = '''0.0849'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').Value = '1.820.99'
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = '1.597.16'
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '''65.22'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '0.0₃0742'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = '''1.01'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''209.44'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = '''7.15'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('D24').Value = '''144.15'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').Value = '''1.01'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '''7.14'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('D29').Value = '''0.0509'
$ws.Range('E29').Value = '  -2.12%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').Value = '''2.99'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('E33').Value = '  +18.13%  '
$ws.Range('D34').Value = '1.279.72'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').Value = '''2.49'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = '''0.597'
$ws.Range('E37').Value = '  -3.82%  '
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('D39').Value = '''0.827'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').Value = '''2.24'
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('D41').Value = '''5.46'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('D42').Value = '''0.779'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').Value = '''62.77'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').Value = '1.733.18'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '''90.55'
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₇0978'
$ws.Range('E49').Value = '  -6.09%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.56'
$ws.Range('E50').Value = '  +2.40%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '''1.01'
$ws.Range('E51').Value = '  +0.43%  '
